$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# Insert a new row above the current row 2 (child_id), shifting everything down.
$ws.Rows.Item(2).EntireRow.Insert()

# Populate the new row 2 with the "row_id" variable definition.
$ws.Range("A2").Value = "row_id"
$ws.Range("B2").Value = "integer"
$ws.Range("C2").Value = "numeric"
$ws.Range("D2").Value = "Unique identifier for the row in Opal"

# child_id (now row 3) valueType changes back from "integer" to "text".
$ws.Range("B3").Value = "text"

# Mirror the saved selection state from the authored workbook.
$ws.Range("A2:D2").Select()
